$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 94, pushing existing row 94 (and everything
# below it) down by one. This matches the diff, where every row from the
# old 94..200 now lives one row lower (95..201), and a brand-new record is
# written into the freshly inserted row 94.
$ws.Rows.Item(94).Insert()

# Populate the new row 94 with the new data record.
$ws.Range("A94").Value = 4
$ws.Range("B94").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C94").Value = "Los Lagos"
$ws.Range("D94").Value = 44539
$ws.Range("E94").Value = 10
$ws.Range("F94").Value = 100112040
$ws.Range("G94").Value = "Cilantro"
$ws.Range("H94").Value = "Sin especificar"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 80
$ws.Range("K94").Value = 6000
$ws.Range("L94").Value = 6000
$ws.Range("M94").Value = 6000
$ws.Range("N94").Value = "$/docena de atados (2 kilos)"
$ws.Range("O94").Value = "Región de La Araucanía"
$ws.Range("P94").Value = 3000
$ws.Range("Q94").Value = 2
$ws.Range("R94").Value = "Hortaliza"
